# Applies: "removed ER tags from non-ER templates and non-ER tags"
#  - rename "SwateTemplateMetadata" sheet to "isa_template"
#  - clear ER value + ER Term Accession Number + ER Term Source REF (B8:B10)
#  - clear Tags Term Source REF values that reference non-tag related terms (D12:D14)
#  - switch the active/selected sheet & selection to the metadata sheet

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("rna_extraction")
$wsMeta = $wb.Worksheets.Item("SwateTemplateMetadata")

# Rename the metadata sheet
$wsMeta.Name = "isa_template"

# Clear the ER row values (B8 = ER, B9 = ER Term Accession Number, B10 = ER Term Source REF)
$wsMeta.Range("B8:B10").ClearContents()

# Clear the non-tag related Term Source REF entries on the Tags rows (D12:D14)
$wsMeta.Range("D12:D14").ClearContents()

# Rows 12 & 13 no longer need the extra wrapped-text height
$wsMeta.Range("A12:A13").EntireRow.AutoFit() | Out-Null

# Update selections / active sheet to match the authored state
$wsData.Activate()
$wsData.Range("W19").Select() | Out-Null
$wsMeta.Activate()
$wsMeta.Range("F10").Select() | Out-Null
